# Search Skll -> Filter, user and any string
# Rebuilds Sheet1 (CategoryIndex/SubCatCount counts table) and Sheet2
# (Category -> SubCategory lookup grid) with the new data set.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------
# Sheet1: rename headers and append the new rows (5,4) (6,6) (7,6) (8,6)
# ---------------------------------------------------------------
$ws1.Range("A1").Value = "CategoryIndex"
$ws1.Range("B1").Value = "SubCatCount"

$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = 4
$ws1.Range("A7").Value = 6
$ws1.Range("B7").Value = 6
$ws1.Range("A8").Value = 7
$ws1.Range("B8").Value = 6
$ws1.Range("A9").Value = 8
$ws1.Range("B9").Value = 6

$ws1.Columns.Item(1).ColumnWidth = 15.666666666666666

# ---------------------------------------------------------------
# Sheet2: Category / SubCategory grid, written column by column so the
# shared-string table is rebuilt in the same order as the target file.
# ---------------------------------------------------------------

# Column A - subcat1
$ws2.Range("A1").Value = "| subcat1                  |"
$ws2.Range("A2").Value = "| Logo Design              |"
$ws2.Range("A3").Value = "| Book Album covers        |"
$ws2.Range("A4").Value = "| Flyers Brochures16       |"
$ws2.Range("A5").Value = "| Web Mobile Design        |"
$ws2.Range("A6").Value = "| Search Display Marketing |"

# Column B - subcat2 (B6 re-uses the same text as A6)
$ws2.Range("B1").Value = "| subcat2                  |"
$ws2.Range("B2").Value = "| Social Media Marketing   |"
$ws2.Range("B3").Value = "| Content Marketing        |"
$ws2.Range("B4").Value = "| Video Marketing          |"
$ws2.Range("B5").Value = "| Email Marketing          |"
$ws2.Range("B6").Value = "| Search Display Marketing |"

# Column C - subcat3
$ws2.Range("C1").Value = "| subcat3               |"
$ws2.Range("C2").Value = "| Resumes Cover Letters |"
$ws2.Range("C3").Value = "| Proof Reading Editing |"
$ws2.Range("C4").Value = "| Translation           |"
$ws2.Range("C5").Value = "| Creative Writing      |"
$ws2.Range("C6").Value = "| Business Copywriting  |"

# Column D - subcat4
$ws2.Range("D1").Value = "| subcat4                 |"
$ws2.Range("D2").Value = "| Promotional Videos      |"
$ws2.Range("D3").Value = "| Editing Post Production |"
$ws2.Range("D4").Value = "| Lyric Music Videos      |"
$ws2.Range("D5").Value = "| Other                   |  "

# Column E - subcat5
$ws2.Range("E1").Value = "| subcat5                |"
$ws2.Range("E2").Value = "| Mixing Mastering       |"
$ws2.Range("E3").Value = "| Voice Over             |"
$ws2.Range("E4").Value = "| Song Writers Composers |"
$ws2.Range("E5").Value = "| Other                  |"

# Column F - subcat6
$ws2.Range("F1").Value = "| subcat6               |"
$ws2.Range("F2").Value = "| WordPress             |"
$ws2.Range("F3").Value = "| Web Mobile App        |"
$ws2.Range("F4").Value = "| Data Analysis Reports |"
$ws2.Range("F5").Value = "| QA                    |"
$ws2.Range("F6").Value = "| Databases             |"
$ws2.Range("F7").Value = "| Other                 |"

# Column G - subcat7
$ws2.Range("G1").Value = "| subcat7              |"
$ws2.Range("G2").Value = "| Business Tips        |"
$ws2.Range("G3").Value = "| Presentations        |"
$ws2.Range("G4").Value = "| Market Advice        |"
$ws2.Range("G5").Value = "| Legal Consulting     |"
$ws2.Range("G6").Value = "| Financial Consulting |"
$ws2.Range("G7").Value = "| Other                |"

# Column H - subcat8
$ws2.Range("H1").Value = "| subcat8                   |"
$ws2.Range("H2").Value = "| Online Lessons            |"
$ws2.Range("H3").Value = "| Relationship Advice       |"
$ws2.Range("H4").Value = "| Astrology                 |"
$ws2.Range("H5").Value = "| Health Nutrition  Fitness |"
$ws2.Range("H6").Value = "| Gaming                    |"
$ws2.Range("H7").Value = "| Other                     |"

# Column widths (closest values this engine's 1/6-character grid can reach)
$ws2.Columns.Item(1).ColumnWidth = 22.5
$ws2.Columns.Item(2).ColumnWidth = 22.5
$ws2.Columns.Item(3).ColumnWidth = 20.5
$ws2.Columns.Item(4).ColumnWidth = 20.833333333333332
$ws2.Columns.Item(5).ColumnWidth = 22.0
$ws2.Columns.Item(6).ColumnWidth = 19.5
$ws2.Columns.Item(7).ColumnWidth = 18.166666666666668
$ws2.Columns.Item(8).ColumnWidth = 21.5

# ---------------------------------------------------------------
# View state: Sheet1 selection B1:B9 (no longer the active tab),
# Sheet2 becomes the active / selected tab with E14 selected and 90% zoom.
# ---------------------------------------------------------------
$ws1.Range("B1:B9").Select() | Out-Null
$ws2.Range("E14").Select() | Out-Null
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 90
